# Auto-generated edit script: refresh market-board derived columns (H..N)
# across all 8 sheets, per scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
# row 11
$ws.Range("H11").Value = 81.22221999999999
$ws.Range("I11").Value = 81.22221999999999
$ws.Range("K11").Value = 81.22221999999999
$ws.Range("M11").Value = 58.77778000000001
# row 34
$ws.Range("H34").Value = 15997.667
$ws.Range("I34").Value = 3997
$ws.Range("J34").Value = 39999
$ws.Range("K34").Value = 3997
$ws.Range("L34").Value = 39999
$ws.Range("M34").Value = -3794
$ws.Range("N34").Value = -40405
# row 36
$ws.Range("H36").Value = 15997.667
$ws.Range("I36").Value = 3997
$ws.Range("J36").Value = 39999
$ws.Range("K36").Value = 3997
$ws.Range("L36").Value = 39999
$ws.Range("M36").Value = -3282
$ws.Range("N36").Value = -41429
# row 97
$ws.Range("H97").Value = 5250.5
$ws.Range("J97").Value = 5250.5
$ws.Range("L97").Value = 15751.5
$ws.Range("N97").Value = -16743.5
# row 111
$ws.Range("H111").Value = 334999.34
$ws.Range("I111").Value = 1499
$ws.Range("J111").Value = 501749.5
$ws.Range("K111").Value = 4497
$ws.Range("L111").Value = 1505248.5
$ws.Range("M111").Value = -1430
$ws.Range("N111").Value = -1511382.5
# row 116
$ws.Range("H116").Value = 5156.645
$ws.Range("I116").Value = 5043
$ws.Range("J116").Value = 5923.75
$ws.Range("K116").Value = 5043
$ws.Range("L116").Value = 5923.75
$ws.Range("M116").Value = -1601
$ws.Range("N116").Value = -12807.75
# row 125
$ws.Range("H125").Value = 19633.166
$ws.Range("J125").Value = 25949.75
$ws.Range("L125").Value = 233547.75
$ws.Range("N125").Value = -238467.75
# row 137
$ws.Range("H137").Value = 2143.5417
$ws.Range("I137").Value = 1826.381
$ws.Range("K137").Value = 5479.143
$ws.Range("M137").Value = -2929.143

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
# row 32
$ws.Range("H32").Value = 3820.7273
$ws.Range("I32").Value = 1827.475
$ws.Range("K32").Value = 1827.475
$ws.Range("M32").Value = -1540.475
# row 45
$ws.Range("H45").Value = 7964
$ws.Range("I45").Value = 9216.333000000001
$ws.Range("K45").Value = 9216.333000000001
$ws.Range("M45").Value = -8839.333000000001
# row 76
$ws.Range("H76").Value = 34607
$ws.Range("J76").Value = 34607
$ws.Range("L76").Value = 34607
$ws.Range("N76").Value = -35283
# row 79
$ws.Range("H79").Value = 34607
$ws.Range("J79").Value = 34607
$ws.Range("L79").Value = 34607
$ws.Range("N79").Value = -36947
# row 80
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
# row 83
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
# row 88
$ws.Range("H88").Value = 3666.5
$ws.Range("I88").Value = 3999.5
$ws.Range("K88").Value = 3999.5
$ws.Range("M88").Value = -3593.5
# row 91
$ws.Range("H91").Value = 3666.5
$ws.Range("I91").Value = 3999.5
$ws.Range("K91").Value = 3999.5
$ws.Range("M91").Value = -2595.5
# row 97
$ws.Range("H97").Value = 885.2
$ws.Range("I97").Value = 882.3333
$ws.Range("K97").Value = 882.3333
$ws.Range("M97").Value = -386.3333

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
# row 86
$ws.Range("H86").Value = 2810.7727
$ws.Range("I86").Value = 2929.9333
$ws.Range("J86").Value = 2555.4285
$ws.Range("K86").Value = 2929.9333
$ws.Range("L86").Value = 2555.4285
$ws.Range("M86").Value = -1806.9333
$ws.Range("N86").Value = -4801.4285
# row 89
$ws.Range("H89").Value = 2810.7727
$ws.Range("I89").Value = 2929.9333
$ws.Range("J89").Value = 2555.4285
$ws.Range("K89").Value = 14649.6665
$ws.Range("L89").Value = 12777.1425
$ws.Range("M89").Value = -9033.666500000001
$ws.Range("N89").Value = -24009.1425
# row 94
$ws.Range("H94").Value = 621
$ws.Range("J94").Value = 499.5
$ws.Range("L94").Value = 499.5
$ws.Range("N94").Value = -1401.5

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
# row 7
$ws.Range("H7").Value = 145.81818
$ws.Range("I7").Value = 187
$ws.Range("K7").Value = 187
$ws.Range("M7").Value = -74
# row 29
$ws.Range("H29").Value = 24999
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 24999
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 24999
$ws.Range("N29").Value = -25585
$ws.Range("M29").ClearContents()
# row 31
$ws.Range("H31").Value = 12402.207
$ws.Range("I31").Value = 10999.6
$ws.Range("K31").Value = 10999.6
$ws.Range("M31").Value = -10704.6
# row 34
$ws.Range("H34").Value = 12402.207
$ws.Range("I34").Value = 10999.6
$ws.Range("K34").Value = 10999.6
$ws.Range("M34").Value = -10797.6
# row 86
$ws.Range("H86").Value = 15071.429
$ws.Range("I86").Value = 12500
$ws.Range("J86").Value = 15500
$ws.Range("K86").Value = 12500
$ws.Range("L86").Value = 15500
$ws.Range("M86").Value = -11377
$ws.Range("N86").Value = -17746
# row 89
$ws.Range("H89").Value = 15071.429
$ws.Range("I89").Value = 12500
$ws.Range("J89").Value = 15500
$ws.Range("K89").Value = 62500
$ws.Range("L89").Value = 77500
$ws.Range("M89").Value = -56884
$ws.Range("N89").Value = -88732
# row 105
$ws.Range("H105").Value = 2042360.8
$ws.Range("I105").Value = 2916686.8
$ws.Range("J105").Value = 2266.6667
$ws.Range("K105").Value = 2916686.8
$ws.Range("L105").Value = 2266.6667
$ws.Range("M105").Value = -2914939.8
$ws.Range("N105").Value = -5760.6667

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
# row 75
$ws.Range("H75").Value = 301.7
$ws.Range("I75").Value = 315.5
$ws.Range("J75").Value = 281
$ws.Range("K75").Value = 946.5
$ws.Range("L75").Value = 843
$ws.Range("M75").Value = 51.5
$ws.Range("N75").Value = -2839
# row 78
$ws.Range("H78").Value = 301.7
$ws.Range("I78").Value = 315.5
$ws.Range("J78").Value = 281
$ws.Range("K78").Value = 2839.5
$ws.Range("L78").Value = 2529
$ws.Range("M78").Value = 2152.5
$ws.Range("N78").Value = -12513
# row 98
$ws.Range("H98").Value = 500.57144
$ws.Range("J98").Value = 424
$ws.Range("L98").Value = 1272
$ws.Range("N98").Value = -4268
# row 99
$ws.Range("H99").Value = 0
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 0
$ws.Range("M99").ClearContents()
$ws.Range("N99").ClearContents()
# row 132
$ws.Range("H132").Value = 21063
$ws.Range("I132").Value = 30344.5
$ws.Range("J132").Value = 2500
$ws.Range("K132").Value = 273100.5
$ws.Range("L132").Value = 22500
$ws.Range("M132").Value = -270570.5
$ws.Range("N132").Value = -27560

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
# row 97
$ws.Range("H97").Value = 1577.2609
$ws.Range("I97").Value = 1362.2
$ws.Range("J97").Value = 1980.5
$ws.Range("K97").Value = 1362.2
$ws.Range("L97").Value = 1980.5
$ws.Range("M97").Value = -866.2
$ws.Range("N97").Value = -2972.5
# row 102
$ws.Range("H102").Value = 2804.6086
$ws.Range("I102").Value = 3008.0527
$ws.Range("J102").Value = 1838.25
$ws.Range("K102").Value = 3008.0527
$ws.Range("L102").Value = 1838.25
$ws.Range("M102").Value = -1386.0527
$ws.Range("N102").Value = -5082.25
# row 113
$ws.Range("H113").Value = 48300.637
$ws.Range("I113").Value = 65369.75
$ws.Range("J113").Value = 2783
$ws.Range("K113").Value = 65369.75
$ws.Range("L113").Value = 2783
$ws.Range("M113").Value = -63199.75
$ws.Range("N113").Value = -7123
# row 126
$ws.Range("H126").Value = 3787.889
$ws.Range("I126").Value = 3532
$ws.Range("J126").Value = 4299.6665
$ws.Range("K126").Value = 10596
$ws.Range("L126").Value = 12898.9995
$ws.Range("M126").Value = -8126
$ws.Range("N126").Value = -17838.9995

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
# row 22
$ws.Range("H22").Value = 3058.6365
$ws.Range("I22").Value = 3114.5
$ws.Range("J22").Value = 2500
$ws.Range("K22").Value = 3114.5
$ws.Range("L22").Value = 2500
$ws.Range("M22").Value = -2819.5
$ws.Range("N22").Value = -3090
# row 27
$ws.Range("H27").Value = 3058.6365
$ws.Range("I27").Value = 3114.5
$ws.Range("J27").Value = 2500
$ws.Range("K27").Value = 3114.5
$ws.Range("L27").Value = 2500
$ws.Range("M27").Value = -3007.5
$ws.Range("N27").Value = -2714
# row 82
$ws.Range("H82").Value = 957.5
$ws.Range("I82").Value = 1090.625
$ws.Range("J82").Value = 744.5
$ws.Range("K82").Value = 1090.625
$ws.Range("L82").Value = 744.5
$ws.Range("M82").Value = -729.625
$ws.Range("N82").Value = -1466.5
# row 85
$ws.Range("H85").Value = 957.5
$ws.Range("I85").Value = 1090.625
$ws.Range("J85").Value = 744.5
$ws.Range("K85").Value = 1090.625
$ws.Range("L85").Value = 744.5
$ws.Range("M85").Value = 157.375
$ws.Range("N85").Value = -3240.5
# row 100
$ws.Range("H100").Value = 10296233
$ws.Range("I100").Value = 35001136
$ws.Range("K100").Value = 35001136
$ws.Range("M100").Value = -35000595

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
# row 100
$ws.Range("H100").Value = 2068.3572
$ws.Range("I100").Value = 2068.3572
$ws.Range("K100").Value = 4136.7144
$ws.Range("M100").Value = -3595.7144
# row 136
$ws.Range("H136").Value = 26317222
$ws.Range("I136").Value = 26317222
$ws.Range("K136").Value = 78951666
$ws.Range("M136").Value = -78949116
